# "Generate Report for Handoff"
# The workbook tracks localization status for a single source file across
# two target locales (zh-cn, de-de), summarized on an "Overview" sheet.
# Re-generating the report after handoff moves the status from
# "In Translation" to "Ready for handoff" and refreshes the handoff
# timestamps. Widen the status/date columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status + the latest handoff xliff date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-13 13:48:48"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-13 13:48:38"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-13 13:48:48"

# --- Widen the columns that now hold the longer "Ready for handoff" text ---
# (ColumnWidth snaps to Excel's pixel grid; 16.3 is the nearest input that
# lands on the widened column size used by the regenerated report.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
